$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.800.28'
$ws.Range('E2').Value = '  -0.20%  '

$ws.Range('D3').Value = '1.684.39'
$ws.Range('E3').Value = '  -1.72%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.67%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '314.07'
$ws.Range('E5').Value = '  -1.41%  '

$ws.Range('E6').Value = '  +0.45%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3934'
$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3971'
$ws.Range('E8').Value = '  -2.52%  '

$ws.Range('B9').Value = 'BinanceUSD'
$ws.Range('C9').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.005'
$ws.Range('E9').Value = '  +0.40%  '

$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.424'
$ws.Range('E10').Value = '  -5.22%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '51.83'
$ws.Range('E11').Value = '  -3.27%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08677'
$ws.Range('E12').Value = '  -1.91%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '25.19'
$ws.Range('E13').Value = '  -5.01%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.322'
$ws.Range('E14').Value = '  -2.85%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.796'
$ws.Range('E15').Value = '  -4.39%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.00001321'
$ws.Range('E16').Value = '  -3.17%  '

$ws.Range('D17').Value = '1.681.71'
$ws.Range('E17').Value = '  -1.70%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '94.19'
$ws.Range('E18').Value = '  -3.32%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.07099'
$ws.Range('E19').Value = '  -1.67%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '20.16'
$ws.Range('E20').Value = '  -2.51%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.140'
$ws.Range('E21').Value = '  -2.63%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.004'

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.10'
$ws.Range('E23').Value = '  -2.45%  '

$ws.Range('D24').Value = '24.781.79'
$ws.Range('E24').Value = '  -0.20%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.369'
$ws.Range('E25').Value = '  +1.44%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '23.85'
$ws.Range('E26').Value = '  +2.25%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.776'
$ws.Range('E27').Value = '  -8.23%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '162.47'
$ws.Range('E28').Value = '  -3.43%  '

$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '150.35'
$ws.Range('E29').Value = '  +2.67%  '

$ws.Range('B30').Value = 'HuobiToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.744'
$ws.Range('E30').Value = '  -3.93%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.597'
$ws.Range('E31').Value = '  +16.79%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.848'
$ws.Range('E32').Value = '  -8.49%  '

$ws.Range('D33').Value = '1.823.27'
$ws.Range('E33').Value = '  -4.00%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.08474'
$ws.Range('E34').Value = '  -4.48%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.03078'
$ws.Range('E35').Value = '  -3.15%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.016'
$ws.Range('E36').Value = '  -4.16%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.949'
$ws.Range('E37').Value = '  -4.21%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2808'
$ws.Range('E38').Value = '  -1.79%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.09577'
$ws.Range('E39').Value = '  +3.44%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '10.49'
$ws.Range('E40').Value = '  -3.85%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.7961'
$ws.Range('E41').Value = '  -6.48%  '

$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.474'
$ws.Range('E42').Value = '  -0.55%  '

$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '13.71'
$ws.Range('E43').Value = '  -3.53%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.80'
$ws.Range('E44').Value = '  -4.62%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.7167'
$ws.Range('E45').Value = '  -4.33%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.586'
$ws.Range('E46').Value = '  -4.58%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.185'
$ws.Range('E47').Value = '  -2.37%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.08767'
$ws.Range('E48').Value = '  +5.49%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.003'

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.344'
$ws.Range('E50').Value = '  -4.75%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '138.71'
$ws.Range('E51').Value = '  -1.41%  '
